$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 136 ("「コーヒー無料です」" post), which shifts all
# subsequent rows up by one.
$ws.Rows.Item(136).Delete()
